$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.098.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.117.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5199"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4458"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09363"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.640"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.121.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.955"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001167"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06694"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.286"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.124.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.322"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.537"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "134.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.787"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1056"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.269"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.638"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.967"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02628"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06861"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.7103"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2249"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.329"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6862"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.91%  "
$ws.Range("E45").Value = "  +4.86%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.268"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.633"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.59%  "
